# Updated symbol list on Thu Dec 22 03:38:17 UTC 2022 with GitHub Actions
#
# Applies the refreshed "Price" (column D) and "Volume(1h)" (column E)
# values scraped for this run. Column D values are numeric-looking text
# (the sheet stores prices as literal strings, e.g. "247.72", not
# numbers), so each one is written with a leading apostrophe to force
# Excel to keep it as text instead of silently re-typing it as a Number
# (which would also strip meaningful trailing zeros like "0.002900").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# Column D (Price) updates
Set-TextValue "D2"  "247.77"
Set-TextValue "D3"  "22.84"
Set-TextValue "D4"  "5.296"
Set-TextValue "D7"  "6.343"
Set-TextValue "D8"  "0.8103"
Set-TextValue "D9"  "0.8753"
Set-TextValue "D10" "0.1433"
Set-TextValue "D11" "0.07412"
Set-TextValue "D13" "0.03113"
Set-TextValue "D14" "0.09387"
Set-TextValue "D15" "3.898"
Set-TextValue "D16" "0.001589"
Set-TextValue "D17" "0.04815"
Set-TextValue "D20" "0.005174"
Set-TextValue "D21" "0.0009972"
Set-TextValue "D23" "3.741"
Set-TextValue "D24" "2.195"
Set-TextValue "D40" "0.03943"
Set-TextValue "D41" "0.006752"
Set-TextValue "D42" "0.1069"
Set-TextValue "D43" "0.002900"
Set-TextValue "D44" "0.008014"
Set-TextValue "D45" "0.00005616"
Set-TextValue "D48" "0.1804"

# Column E (Volume(1h)) updates
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Write-Host "Applied cryptos.xlsx price/volume refresh"
